# re-render with lisa edits to ch1 and 2
# Shift the date values in column F (rows 2-7) forward by 10 days,
# keeping their existing date number-format/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 44507
$ws.Range("F3").Value = 44506
$ws.Range("F4").Value = 44505
$ws.Range("F5").Value = 44504
$ws.Range("F6").Value = 44503
$ws.Range("F7").Value = 44502
